# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (row 1) ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 22 de Marzo de 2020 a las 18:16"

# --- Update numeric data for existing countries ---

# Row 5: Italia
$ws.Cells.Item(5,2).Value = 59138
$ws.Cells.Item(5,3).Value = 5560
$ws.Cells.Item(5,4).Value = 7024
$ws.Cells.Item(5,5).Value = 46638
$ws.Cells.Item(5,6).Value = 3000
$ws.Cells.Item(5,7).Value = 651
$ws.Cells.Item(5,8).Value = 5476

# Row 6: Estados Unidos
$ws.Cells.Item(6,2).Value = 30291
$ws.Cells.Item(6,3).Value = 6084
$ws.Cells.Item(6,5).Value = 29725

# Row 8: Alemania
$ws.Cells.Item(8,2).Value = 24714
$ws.Cells.Item(8,3).Value = 2350
$ws.Cells.Item(8,5).Value = 24356
$ws.Cells.Item(8,7).Value = 8
$ws.Cells.Item(8,8).Value = 92

# Row 13: Reino Unido
$ws.Cells.Item(13,2).Value = 5683
$ws.Cells.Item(13,3).Value = 665
$ws.Cells.Item(13,5).Value = 5309
$ws.Cells.Item(13,7).Value = 48
$ws.Cells.Item(13,8).Value = 281

# Row 17: Noruega
$ws.Cells.Item(17,2).Value = 2263
$ws.Cells.Item(17,3).Value = 99
$ws.Cells.Item(17,5).Value = 2250

# --- Re-order "Serbia" to sit right after "Argentina" (row 62), pushing
# "Republica Dominicana", "Argelia" and "Armenia" down by one row each,
# matching the new shared-string order: Argentina, Serbia, Republica
# Dominicana, Argelia, Armenia, Kuwait. Row 61 (Argentina) and row 66
# (Kuwait) are untouched; only the data in rows 62-65 is rewritten.

# Row 65 <- old row 64 (Argelia)
$ws.Cells.Item(65,1).Value = "Armenia"
$ws.Cells.Item(65,2).Value = 190
$ws.Cells.Item(65,3).Value = 30
$ws.Cells.Item(65,4).Value = 2
$ws.Cells.Item(65,5).Value = 188
$ws.Cells.Item(65,6).Value = 6
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(65,8).Value = 0

# Row 64 <- old row 63 (Republica Dominicana)
$ws.Cells.Item(64,1).Value = "Argelia"
$ws.Cells.Item(64,2).Value = 201
$ws.Cells.Item(64,3).Value = 62
$ws.Cells.Item(64,4).Value = 65
$ws.Cells.Item(64,5).Value = 119
$ws.Cells.Item(64,6).Value = 0
$ws.Cells.Item(64,7).Value = 2
$ws.Cells.Item(64,8).Value = 17

# Row 63 <- old row 62 (Sudafrica's slot, here was "Republica Dominicana")
$ws.Cells.Item(63,1).Value = "Republica Dominicana"
$ws.Cells.Item(63,2).Value = 202
$ws.Cells.Item(63,3).Value = 90
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 199
$ws.Cells.Item(63,6).Value = 0
$ws.Cells.Item(63,7).Value = 0
$ws.Cells.Item(63,8).Value = 3

# Row 62 <- new "Serbia" row
$ws.Cells.Item(62,1).Value = "Serbia"
$ws.Cells.Item(62,2).Value = 222
$ws.Cells.Item(62,3).Value = 51
$ws.Cells.Item(62,4).Value = 2
$ws.Cells.Item(62,5).Value = 218
$ws.Cells.Item(62,6).Value = 4
$ws.Cells.Item(62,7).Value = 1
$ws.Cells.Item(62,8).Value = 2

$wb.Save()
